$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.300.04"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "1.929.30"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'249.51"
$ws.Range("E5").Value = "  -1.20%  "
$ws.Range("D6").Value = "'0.7202"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("D7").Value = "'0.9997"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'27.99"
$ws.Range("E8").Value = "  -2.36%  "
$ws.Range("D9").Value = "'0.3204"
$ws.Range("E9").Value = "  -4.11%  "
$ws.Range("D10").Value = "'0.07099"
$ws.Range("E10").Value = "  -3.49%  "
$ws.Range("D11").Value = "'0.7885"
$ws.Range("E11").Value = "  -3.32%  "
$ws.Range("D12").Value = "'0.08010"
$ws.Range("E12").Value = "  -1.75%  "
$ws.Range("D13").Value = "1.932.80"
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").Value = "'5.377"
$ws.Range("E14").Value = "  -2.10%  "
$ws.Range("D15").Value = "'94.71"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("D16").Value = "'14.66"
$ws.Range("E16").Value = "  -1.53%  "
$ws.Range("D17").Value = "30.298.84"
$ws.Range("D18").Value = "'257.03"
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("D19").Value = "'0.000008071"
$ws.Range("E19").Value = "  -3.65%  "
$ws.Range("D20").Value = "'5.732"
$ws.Range("E20").Value = "  -2.13%  "
$ws.Range("D21").Value = "2.181.35"
$ws.Range("E21").Value = "  -0.68%  "
$ws.Range("D22").Value = "'0.9994"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "'6.816"
$ws.Range("E24").Value = "  -2.27%  "
$ws.Range("D25").Value = "'9.565"
$ws.Range("E25").Value = "  -2.92%  "
$ws.Range("D26").Value = "'164.58"
$ws.Range("E26").Value = "  +1.87%  "
$ws.Range("D27").Value = "'19.09"
$ws.Range("E27").Value = "  -1.57%  "
$ws.Range("E28").Value = "  -4.87%  "
$ws.Range("D29").Value = "'0.1282"
$ws.Range("E29").Value = "  -2.50%  "
$ws.Range("D30").Value = "'1.360"
$ws.Range("E31").Value = "  -2.84%  "
$ws.Range("D32").Value = "'4.418"
$ws.Range("E32").Value = "  -1.57%  "
$ws.Range("D33").Value = "'4.147"
$ws.Range("E33").Value = "  -2.83%  "
$ws.Range("D34").Value = "'0.05117"
$ws.Range("E34").Value = "  -2.90%  "
$ws.Range("D35").Value = "'1.287"
$ws.Range("E35").Value = "  +0.76%  "
$ws.Range("D36").Value = "'0.7488"
$ws.Range("E36").Value = "  -1.77%  "
$ws.Range("D37").Value = "'2.768"
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("D38").Value = "'0.01987"
$ws.Range("D39").Value = "'2.797"
$ws.Range("E39").Value = "  -1.78%  "
$ws.Range("D40").Value = "'78.28"
$ws.Range("E40").Value = "  -3.76%  "
$ws.Range("D41").Value = "'6.399"
$ws.Range("E41").Value = "  -2.56%  "
$ws.Range("D42").Value = "'0.4524"
$ws.Range("E42").Value = "  -0.97%  "
$ws.Range("D43").Value = "'1.994"
$ws.Range("E43").Value = "  -2.03%  "
$ws.Range("D44").Value = "'0.8460"
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").Value = "'101.07"
$ws.Range("E46").Value = "  -2.09%  "
$ws.Range("D47").Value = "'9.835"
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("D48").Value = "'7.487"
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("D49").Value = "'36.86"
$ws.Range("E49").Value = "  -0.72%  "
$ws.Range("D50").Value = "'965.47"
$ws.Range("E50").Value = "  +9.36%  "
$ws.Range("D51").Value = "'0.4218"
$ws.Range("E51").Value = "  +0.67%  "
